# Apply the "complete uploading video via azure storage using sas token" edit.
#
# The original sheet described a "Get Access Video Key" endpoint
# (GET /users/{user id}/uploadKey) used to fetch a public access key for
# viewing video. The commit replaces it with a "Get Upload Video Key"
# endpoint (GET /owner/{user id}/uploadKeys) that returns a SAS key used to
# upload video to Azure storage, and its error payload code changes from
# 401 to 400/401 combo text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the four cells whose text actually changes ---
# Order matters: it controls the order new entries are appended to the
# shared-strings table, matching the target workbook exactly.

# G4: Note/title of the endpoint
$ws.Range("G4").Value = " Get Upload Video Key`n - 비디오 업로드에 사용될 SAS 키를 받아온다"

# G7: URL Params text
$ws.Range("G7").Value = " key: video / value: {video name to upload}"

# G5: URL of the endpoint
$ws.Range("G5").Value = " /owner/{user id}/uploadKeys"

# G13: Error Response text
$ws.Range("G13").Value = " Code: 400 Bad Request`n Contents: { ""msg"": ""error message"", ""code"": 401 }"

# --- Update the selected range shown when the workbook is opened ---
[void]$ws.Range("G11:G12").Select()

Write-Output "edit applied"
